$wb = $excel.ActiveWorkbook

# Portfolio sheet - convert A2:A10 from text to numeric values
$wsPortfolio = $wb.Worksheets.Item("Portfolio")
$tickers = @(394670, 292150, 483420, 245350, 469160, 220130, 419430, 105010, 455960)
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2
    $wsPortfolio.Cells.Item($row, 1).Value = $tickers[$i]
}

# History sheet - update the date text in A2 (kept as literal text, not
# auto-converted to a date serial number by Excel's smart-entry parsing)
$wsHistory = $wb.Worksheets.Item("History")
$wsHistory.Range("A2").NumberFormat = "@"
$wsHistory.Range("A2").Value = "2026-02-13"
$wsHistory.Range("A2").Style = "Normal"
